# "Tweaked report values, reset total difference"
#
# The bullet list had an accidental duplicate pair of items ("Better
# feedback" / "Ensure correct conversion") sitting right above the real
# "Polish" / "Make the page look nicer" pair. Removing the duplicate pair
# lets every bullet below it shift up into the slot vacated by its twin.
#
# Separately, the trailing "Separate score tables for game lengths" bullet
# is split into two runs (break after "Separate score tables ") and the
# document's `_GoBack` bookmark - previously sitting at the very end of the
# "Multiplayer?" bullet - is relocated to that split point.

$d = $word.ActiveDocument

# --- Drop the duplicate "Better feedback" / "Ensure correct conversion" pair ---
# (delete the later one first so the earlier index stays valid)
$d.Paragraphs.Item(6).Range.Delete()   # "Ensure correct conversion"
$d.Paragraphs.Item(5).Range.Delete()   # "Better feedback"

# --- Split the "Separate score tables for game lengths" bullet and move the bookmark ---
$scoreTablesPara = $d.Paragraphs.Item(12)
$prefix = "Separate score tables "
$splitPos = $scoreTablesPara.Range.Start + $prefix.Length

$d.Bookmarks.Item("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))
